{"js": "// \"Added button to the header\" - remove the extra contact-info lines\n// (address, phone, email) from the resume header block, leaving just\n// the name and the (now left-aligned) blank spacer paragraph above\n// \"Professional Summary\".\n\n// Locate the header paragraphs by their known text content so the\n// script is resilient to exact paragraph indices.\nconst addressResults = context.document.body.search(\"Lehi, UT 84043\", { matchCase: false });\naddressResults.load(\"items\");\nawait context.sync();\n\nif (addressResults.items.length === 0) {\n  return \"address line not found\";\n}\n\nconst addressPara = addressResults.items[0].paragraphs.getFirst();\naddressPara.load(\"text\");\nawait context.sync();\n\n// The phone and email lines are the two paragraphs immediately\n// following the address line, and the (already blank) centered\n// spacer paragraph follows those.\nconst phonePara = addressPara.getNext();\nconst emailPara = phonePara.getNext();\nconst spacerPara = emailPara.getNext();\n\nphonePara.load(\"text\");\nemailPara.load(\"text\");\nspacerPara.load(\"text,alignment\");\nawait context.sync();\n\n// Delete the address, phone and email paragraphs entirely.\naddressPara.delete();\nphonePara.delete();\nemailPara.delete();\n\n// The spacer paragraph (formerly centered to match the contact block)\n// reverts to left alignment now that it directly follows the name.\nspacerPara.alignment = Word.Alignment.left;\n\nawait context.sync();\n", "ps1": "# \"Added button to the header\" - remove the extra contact-info lines\n# (address, phone, email) from the resume header block, leaving just\n# the name and the (now left-aligned) blank spacer paragraph above\n# \"Professional Summary\".\n\n$d = $word.ActiveDocument\n\n# Find the \"Lehi, UT 84043\" address paragraph by its text so the script\n# does not depend on a hard-coded paragraph index.\n$addressIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"Lehi, UT 84043*\") {\n        $addressIndex = $i\n        break\n    }\n}\n\nif ($addressIndex -eq -1) {\n    throw \"Could not find the address paragraph (Lehi, UT 84043)\"\n}\n\n# The phone number and email paragraphs immediately follow the address\n# line, and the blank, centered spacer paragraph follows those.\n#   $addressIndex     -> \"Lehi, UT 84043\"\n#   $addressIndex + 1 -> \"(801) 369-3018 \"\n#   $addressIndex + 2 -> \"howeti@gmail.com\"\n#   $addressIndex + 3 -> \"\" (divdocumentdivaddressbottomdiv spacer, centered)\n\n# Deleting a paragraph's Range shifts every later paragraph up by one, so\n# re-fetching the same index after each delete walks address -> phone -> email.\n$d.Paragraphs.Item($addressIndex).Range.Delete()\n$d.Paragraphs.Item($addressIndex).Range.Delete()\n$d.Paragraphs.Item($addressIndex).Range.Delete()\n\n# The spacer paragraph (formerly centered to match the contact block)\n# reverts to left alignment now that it directly follows the name,\n# which removes the now-redundant <w:jc w:val=\"center\"/>.\n$spacer = $d.Paragraphs.Item($addressIndex)\n$spacer.Alignment = 0\n"}
